# edit.ps1 - Applies the "Correccion de la respuesta en el punto 1" change
#
# Summary of changes:
#  1. Merge the split runs of the question paragraph
#     "Que significa el contador de referencias..." into a single run (no text change).
#  2. Rewrite the answer paragraph "R/ El contador de referencias..." so that:
#       - "R/" becomes its own run
#       - the _GoBack bookmark now starts right after "R/"
#       - the explanation text is corrected/expanded
#       - "hards" is fixed to "hard"
#       - the _GoBack bookmark ends at the very end of the paragraph
#  3. Move the stray _GoBack bookmark (previously sitting in an empty
#     paragraph near the end of the document) away from there, since it now
#     lives in the answer paragraph above.
#  4. Merge a few other split runs (pure run-merges, no textual change):
#       - "Creamos una funcion..." paragraph
#       - "Cuando se invoca con la bandera -l..." paragraph
#       - "Y ya por ultimo ponemos la bandera -l..." paragraph (keeping the
#         leading <w:lastRenderedPageBreak/>)
#       - "Este programa requiere de 3 elementos..." paragraph

$d = $word.ActiveDocument

function Set-ParaInnerXml {
    param(
        [__ComObject]$Para,
        [string]$InnerXml
    )
    $pStart = $Para.Range.Start
    $pEnd   = $Para.Range.End
    # Trim off the trailing paragraph mark so only the paragraph's content
    # (runs, proofErr markers, bookmarks, ...) gets replaced.
    $r = $d.Range($pStart, $pEnd - 1)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $InnerXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1. Question paragraph: "Que significa el contador de referencias..."
#    Just merge the two runs into a single run with identical text.
# ---------------------------------------------------------------------------
$pQuestion = $d.Paragraphs.Item(6)
Set-ParaInnerXml $pQuestion ('<w:r><w:t>' +
    '¿Qué significa el contador de referencias de un directorio?, ¿cómo cambia el número de entradas en el directorio?' +
    '</w:t></w:r>')

# ---------------------------------------------------------------------------
# 2. Remove the old (stray) _GoBack bookmark before re-adding it below, so
#    there is never more than one bookmark with that name at once.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 3. Answer paragraph: "R/ El contador de referencias..."
# ---------------------------------------------------------------------------
$pAnswer = $d.Paragraphs.Item(7)
$answerXml =
    '<w:r><w:t>R/</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:r><w:t xml:space="preserve"> El contador de referencias se refiere a las </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Hard</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> links que este tiene, y mientras más </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">apuntadores a esa carpeta se tengan, más número de </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>hard</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> links existirán</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:bookmarkEnd w:id="0"/>'
Set-ParaInnerXml $pAnswer $answerXml

# ---------------------------------------------------------------------------
# 4. "Creamos una funcion que lista los archivos..." - merge the two runs.
# ---------------------------------------------------------------------------
$pCreamos = $d.Paragraphs.Item(24)
Set-ParaInnerXml $pCreamos ('<w:r><w:t>' +
    'Creamos una función que lista los archivos de un directorio dado. Cuando es llamado sin argumentos como en este caso, el programa simplemente imprime solo los nombres de los archivos.' +
    '</w:t></w:r>')

# ---------------------------------------------------------------------------
# 5. "Cuando se invoca con la bandera -l..." - merge the first two runs,
#    keep the proofErr-wrapped "stat(" runs untouched.
# ---------------------------------------------------------------------------
$pCuando = $d.Paragraphs.Item(26)
$cuandoXml =
    '<w:r><w:t xml:space="preserve">Cuando se invoca con la bandera -l el programa dará la información de cada uno de los archivos de la carpeta en la que estamos. La información de cada archivo se consigue con el método </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>stat</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>(</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>) que creamos anteriormente.</w:t></w:r>'
Set-ParaInnerXml $pCuando $cuandoXml

# ---------------------------------------------------------------------------
# 6. "Y ya por ultimo ponemos la bandera -l..." - merge the two runs while
#    keeping the leading <w:lastRenderedPageBreak/>.
# ---------------------------------------------------------------------------
$pYaUltimo = $d.Paragraphs.Item(28)
Set-ParaInnerXml $pYaUltimo ('<w:r><w:lastRenderedPageBreak/><w:t>' +
    'Y ya por último ponemos la bandera -l y una ruta para saber la información de los archivos de la ruta que le pondremos.' +
    '</w:t></w:r>')

# ---------------------------------------------------------------------------
# 7. "Este programa requiere de 3 elementos..." - merge the two runs.
# ---------------------------------------------------------------------------
$pEste = $d.Paragraphs.Item(51)
$quote1 = [char]0x201C
$quote2 = [char]0x201D
Set-ParaInnerXml $pEste ('<w:r><w:t>' +
    'Este programa requiere de 3 elementos a la hora de la ejecución, necesitamos el nombre del archivo a ejecutar, el número de líneas que quieres del final del archivo, acompañadas de un ' +
    $quote1 + '-' + $quote2 +
    ' (-3) y la dirección del archivo. Y así imprimimos en este caso las 10 últimas líneas del mismísimo código.' +
    '</w:t></w:r>')

Write-Host "Edit complete"
